# Updated legacy GSC export data: drop the two oldest days
# (2025-09-19 and 2025-09-20) from the "Chart" data table. Deleting the
# entire rows shifts every following date row up by two, shrinking the
# used range from A1:C89 to A1:C87 and dropping the two now-unused date
# strings from the shared-string table (their removal also renumbers the
# string indices used on the "Table" sheet, which is handled automatically
# since that sheet still references the same text values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 = 2025-09-19, Row 3 = 2025-09-20 (the two oldest entries).
# Deleting both entire rows shifts all subsequent rows up by two.
$ws.Range("A2:A3").EntireRow.Delete()
